$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: checkedIn row - booleans instead of hello/goodbye placeholder text
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = $false

# Row 5: verified row - new boolean cells
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = $false

# Row 6: attendanceNote row - new note text cells (added before row 3's date
# strings so shared-string ordering matches the target)
$ws.Range("B6").Value = "My name is Jonas"
$ws.Range("C6").Value = "I'm carrying the will"

# Row 3: attendanceDate row - add quote-prefixed, date-like text values
$ws.Range("B3").Value = "'02/03/2017"
$ws.Range("C3").Value = "'01/22/2017"
$ws.Range("C3").NumberFormat = "mm-dd-yy"

# Column C width adjustment
$ws.Columns.Item(3).ColumnWidth = 16.5

# Update selection/view (also clears the stale topLeftCell scroll position)
$ws.Range("C15").Select() | Out-Null
